$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.246684350132626
$ws.Cells.Item(2, 3).Value = 0.4748010610079575
$ws.Cells.Item(2, 10).Value = 0.01856763925729443
$ws.Cells.Item(2, 16).Value = 0.1803713527851459
$ws.Cells.Item(2, 19).Value = 0.07957559681697612
$ws.Cells.Item(3, 2).Value = 0.02590673575129534
$ws.Cells.Item(3, 3).Value = 0.07253886010362694
$ws.Cells.Item(3, 10).Value = 0.04663212435233161
$ws.Cells.Item(3, 16).Value = 0.7098445595854922
$ws.Cells.Item(3, 19).Value = 0.1450777202072539
$ws.Cells.Item(4, 10).Value = 0.04651162790697674
$ws.Cells.Item(4, 16).Value = 0.6046511627906976
$ws.Cells.Item(4, 19).Value = 0.3488372093023256
$ws.Cells.Item(6, 2).Value = 0.08205128205128205
$ws.Cells.Item(6, 4).Value = 0.01025641025641026
$ws.Cells.Item(6, 6).Value = 0.03076923076923077
$ws.Cells.Item(6, 10).Value = 0.241025641025641
$ws.Cells.Item(6, 15).Value = 0.04102564102564103
$ws.Cells.Item(6, 17).Value = 0.1230769230769231
$ws.Cells.Item(6, 18).Value = 0.07179487179487179
$ws.Cells.Item(6, 19).Value = 0.4
$ws.Cells.Item(7, 2).Value = 0.1005291005291005
$ws.Cells.Item(7, 4).Value = 0.02645502645502645
$ws.Cells.Item(7, 6).Value = 0.0582010582010582
$ws.Cells.Item(7, 10).Value = 0.1428571428571428
$ws.Cells.Item(7, 15).Value = 0.02116402116402116
$ws.Cells.Item(7, 17).Value = 0.1428571428571428
$ws.Cells.Item(7, 18).Value = 0.06878306878306878
$ws.Cells.Item(7, 19).Value = 0.4391534391534391
$ws.Cells.Item(8, 2).Value = 0.0718562874251497
$ws.Cells.Item(8, 4).Value = 0.02794411177644711
$ws.Cells.Item(8, 5).Value = 0.001996007984031936
$ws.Cells.Item(8, 6).Value = 0.04191616766467066
$ws.Cells.Item(8, 10).Value = 0.12375249500998
$ws.Cells.Item(8, 15).Value = 0.02994011976047904
$ws.Cells.Item(8, 17).Value = 0.1836327345309381
$ws.Cells.Item(8, 18).Value = 0.07984031936127745
$ws.Cells.Item(8, 19).Value = 0.4391217564870259
$ws.Cells.Item(9, 2).Value = 0.08370044052863436
$ws.Cells.Item(9, 4).Value = 0.02202643171806168
$ws.Cells.Item(9, 6).Value = 0.05286343612334802
$ws.Cells.Item(9, 10).Value = 0.1101321585903084
$ws.Cells.Item(9, 15).Value = 0.02202643171806168
$ws.Cells.Item(9, 17).Value = 0.1894273127753304
$ws.Cells.Item(9, 18).Value = 0.08370044052863436
$ws.Cells.Item(9, 19).Value = 0.4361233480176211
$ws.Cells.Item(10, 2).Value = 0.1231884057971015
$ws.Cells.Item(10, 4).Value = 0.01383399209486166
$ws.Cells.Item(10, 6).Value = 0.05072463768115942
$ws.Cells.Item(10, 10).Value = 0.1429512516469038
$ws.Cells.Item(10, 15).Value = 0.02108036890645586
$ws.Cells.Item(10, 17).Value = 0.2015810276679842
$ws.Cells.Item(10, 18).Value = 0.07312252964426877
$ws.Cells.Item(10, 19).Value = 0.3735177865612648
$ws.Cells.Item(11, 7).Value = 0.1634615384615385
$ws.Cells.Item(11, 10).Value = 0.1185897435897436
$ws.Cells.Item(11, 11).Value = 0.2307692307692308
$ws.Cells.Item(11, 12).Value = 0.4743589743589743
$ws.Cells.Item(11, 19).Value = 0.01282051282051282
$ws.Cells.Item(12, 7).Value = 0.7098765432098766
$ws.Cells.Item(12, 10).Value = 0.191358024691358
$ws.Cells.Item(12, 11).Value = 0.006172839506172839
$ws.Cells.Item(12, 12).Value = 0.06790123456790123
$ws.Cells.Item(12, 19).Value = 0.02469135802469136
$ws.Cells.Item(13, 7).Value = 0.5283018867924528
$ws.Cells.Item(13, 10).Value = 0.4339622641509434
$ws.Cells.Item(13, 19).Value = 0.03773584905660377
$ws.Cells.Item(15, 6).Value = 0.02183406113537118
$ws.Cells.Item(15, 8).Value = 0.1703056768558952
$ws.Cells.Item(15, 9).Value = 0.05240174672489083
$ws.Cells.Item(15, 10).Value = 0.3406113537117904
$ws.Cells.Item(15, 11).Value = 0.03056768558951965
$ws.Cells.Item(15, 13).Value = 0.01746724890829694
$ws.Cells.Item(15, 15).Value = 0.07423580786026202
$ws.Cells.Item(15, 19).Value = 0.2925764192139738
$ws.Cells.Item(16, 6).Value = 0.02727272727272727
$ws.Cells.Item(16, 8).Value = 0.1727272727272727
$ws.Cells.Item(16, 9).Value = 0.07727272727272727
$ws.Cells.Item(16, 10).Value = 0.4181818181818182
$ws.Cells.Item(16, 11).Value = 0.1045454545454545
$ws.Cells.Item(16, 13).Value = 0.01363636363636364
$ws.Cells.Item(16, 15).Value = 0.03636363636363636
$ws.Cells.Item(16, 19).Value = 0.15
$ws.Cells.Item(17, 6).Value = 0.01414141414141414
$ws.Cells.Item(17, 8).Value = 0.1838383838383838
$ws.Cells.Item(17, 9).Value = 0.1111111111111111
$ws.Cells.Item(17, 10).Value = 0.4585858585858586
$ws.Cells.Item(17, 11).Value = 0.08888888888888889
$ws.Cells.Item(17, 13).Value = 0.0202020202020202
$ws.Cells.Item(17, 14).Value = 0.00202020202020202
$ws.Cells.Item(17, 15).Value = 0.05454545454545454
$ws.Cells.Item(17, 19).Value = 0.06666666666666667
$ws.Cells.Item(18, 6).Value = 0.005154639175257732
$ws.Cells.Item(18, 8).Value = 0.2216494845360825
$ws.Cells.Item(18, 9).Value = 0.06185567010309279
$ws.Cells.Item(18, 10).Value = 0.4381443298969072
$ws.Cells.Item(18, 11).Value = 0.1185567010309278
$ws.Cells.Item(18, 13).Value = 0.02061855670103093
$ws.Cells.Item(18, 15).Value = 0.08247422680412371
$ws.Cells.Item(18, 19).Value = 0.05154639175257732
$ws.Cells.Item(19, 6).Value = 0.0130718954248366
$ws.Cells.Item(19, 8).Value = 0.2127814088598402
$ws.Cells.Item(19, 9).Value = 0.09658678286129267
$ws.Cells.Item(19, 10).Value = 0.4132171387073348
$ws.Cells.Item(19, 11).Value = 0.1016702977487291
$ws.Cells.Item(19, 13).Value = 0.02396514161220044
$ws.Cells.Item(19, 14).Value = 0.0007262164124909223
$ws.Cells.Item(19, 15).Value = 0.05301379811183732
$ws.Cells.Item(19, 19).Value = 0.08496732026143791
